# Applies the "cryptos list" data refresh (GitHub Actions bot update).
# Updates Price (D) and Volume(1h) (E) columns for most rows, and for
# rows 48-49 additionally swaps the Coin/Link (B/C) values (EOS <-> Quant).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.082.67"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.02%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.875.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.33%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.21%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.20%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.14%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5047"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.36%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3961"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.38%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08210"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.72%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.12"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.70%  "

# Row 11
$ws.Range("E11").Value = "  -3.02%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.57"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.08%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.869.59"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.90%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.290"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.25%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.194"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.00%  "

# Row 16
$ws.Range("E16").Value = "  +0.08%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.94"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.86%  "

# Row 18
$ws.Range("E18").Value = "  -2.44%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.46%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.20%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "30.072.04"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.09%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.828"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.17%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.13"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.85%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.167"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.92%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.085.05"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.01%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.28%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "160.60"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.30%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.238"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -8.76%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.35"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.62%  "

# Row 31
$ws.Range("E31").Value = "  -0.04%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1035"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.97%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.942"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.04%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.703"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.80%  "

# Row 35
$ws.Range("E35").Value = "  -2.45%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.278"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.63%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06372"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.88%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2135"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.76%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.171"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.39%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.508"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.50%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6297"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.92%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.215"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.53%  "

# Row 43
$ws.Range("E43").Value = "  -3.26%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.14"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.40%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5907"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.18%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.092"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.34%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.629"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.81%  "

# Row 48
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "122.17"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.02%  "

# Row 49
$ws.Range("B49").Value = "EOS"
$ws.Range("C49").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.207"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.34%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "77.52"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.87%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.117"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.25%  "
